$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new hours values for D29 and D32 (part of reorganizing "Eval Client" entries)
$ws.Range("D29").Value = 5
$ws.Range("D32").Value = 5

# Update view state: scroll so row 13 is the top-left visible cell, and move the
# active selection to D27
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("D27").Select()
